$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 66.47695399999999
$ws.Cells.Item(2, 8).Value = 199.430862
$ws.Cells.Item(2, 9).Value = 0.04311983106164722
$ws.Cells.Item(2, 10).Value = 0.04311983106164721
$ws.Cells.Item(2, 13).Value = 19.21315233333334
$ws.Cells.Item(2, 14).Value = 57.63945700000001
$ws.Cells.Item(2, 15).Value = 0.04451179209991234
$ws.Cells.Item(2, 16).Value = 0.04451179209991233
$ws.Cells.Item(2, 17).Value = 1277.231843857993
$ws.Cells.Item(2, 18).Value = 11495.08659472194
$ws.Cells.Item(2, 19).Value = 0.001919340955599383
$ws.Cells.Item(2, 20).Value = 0.001919340955599383

$ws.Cells.Item(3, 7).Value = 66.47695399999999
$ws.Cells.Item(3, 8).Value = 199.430862
$ws.Cells.Item(3, 9).Value = 0.04311983106164722
$ws.Cells.Item(3, 10).Value = 0.04311983106164721
$ws.Cells.Item(3, 15).Value = 0.2141755495962477
$ws.Cells.Item(3, 16).Value = 0.2141755495962477
$ws.Cells.Item(3, 17).Value = 6145.603652759989
$ws.Cells.Item(3, 18).Value = 55310.4328748399
$ws.Cells.Item(3, 19).Value = 0.009235213516125648
$ws.Cells.Item(3, 20).Value = 0.009235213516125645

$ws.Cells.Item(4, 7).Value = 66.47695399999999
$ws.Cells.Item(4, 8).Value = 199.430862
$ws.Cells.Item(4, 9).Value = 0.04311983106164722
$ws.Cells.Item(4, 10).Value = 0.04311983106164721
$ws.Cells.Item(4, 13).Value = 166.8580016666666
$ws.Cells.Item(4, 14).Value = 500.5740049999999
$ws.Cells.Item(4, 15).Value = 0.3865658561145097
$ws.Cells.Item(4, 16).Value = 0.3865658561145097
$ws.Cells.Item(4, 17).Value = 11092.21170132692
$ws.Cells.Item(4, 18).Value = 99829.90531194229
$ws.Cells.Item(4, 19).Value = 0.01666865440985869
$ws.Cells.Item(4, 20).Value = 0.01666865440985868

$ws.Cells.Item(5, 7).Value = 66.47695399999999
$ws.Cells.Item(5, 8).Value = 199.430862
$ws.Cells.Item(5, 9).Value = 0.04311983106164722
$ws.Cells.Item(5, 10).Value = 0.04311983106164721
$ws.Cells.Item(5, 13).Value = 41.09915599999999
$ws.Cells.Item(5, 14).Value = 123.297468
$ws.Cells.Item(5, 15).Value = 0.09521587377309249
$ws.Cells.Item(5, 16).Value = 0.09521587377309249
$ws.Cells.Item(5, 17).Value = 2732.146702850823
$ws.Cells.Item(5, 18).Value = 24589.32032565741
$ws.Cells.Item(5, 19).Value = 0.004105692391482874
$ws.Cells.Item(5, 20).Value = 0.004105692391482873

$ws.Cells.Item(6, 7).Value = 66.47695399999999
$ws.Cells.Item(6, 8).Value = 199.430862
$ws.Cells.Item(6, 9).Value = 0.04311983106164722
$ws.Cells.Item(6, 10).Value = 0.04311983106164721
$ws.Cells.Item(6, 13).Value = 112.0244103333333
$ws.Cells.Item(6, 14).Value = 336.073231
$ws.Cells.Item(6, 15).Value = 0.2595309284162377
$ws.Cells.Item(6, 16).Value = 0.2595309284162377
$ws.Cells.Item(6, 17).Value = 7447.041572606124
$ws.Cells.Item(6, 18).Value = 67023.37415345512
$ws.Cells.Item(6, 19).Value = 0.01119092978858063
$ws.Cells.Item(6, 20).Value = 0.01119092978858063

$ws.Cells.Item(7, 9).Value = 0.8830494168872806
$ws.Cells.Item(7, 10).Value = 0.8830494168872804
$ws.Cells.Item(7, 13).Value = 19.21315233333334
$ws.Cells.Item(7, 14).Value = 57.63945700000001
$ws.Cells.Item(7, 15).Value = 0.04451179209991234
$ws.Cells.Item(7, 16).Value = 0.04451179209991233
$ws.Cells.Item(7, 17).Value = 26156.38343610851
$ws.Cells.Item(7, 18).Value = 235407.4509249766
$ws.Cells.Item(7, 19).Value = 0.03930611205843545
$ws.Cells.Item(7, 20).Value = 0.03930611205843544

$ws.Cells.Item(8, 9).Value = 0.8830494168872806
$ws.Cells.Item(8, 10).Value = 0.8830494168872804
$ws.Cells.Item(8, 15).Value = 0.2141755495962477
$ws.Cells.Item(8, 16).Value = 0.2141755495962477
$ws.Cells.Item(8, 19).Value = 0.1891275941824794
$ws.Cells.Item(8, 20).Value = 0.1891275941824793

$ws.Cells.Item(9, 9).Value = 0.8830494168872806
$ws.Cells.Item(9, 10).Value = 0.8830494168872804
$ws.Cells.Item(9, 13).Value = 166.8580016666666
$ws.Cells.Item(9, 14).Value = 500.5740049999999
$ws.Cells.Item(9, 15).Value = 0.3865658561145097
$ws.Cells.Item(9, 16).Value = 0.3865658561145097
$ws.Cells.Item(9, 17).Value = 227156.9909641671
$ws.Cells.Item(9, 18).Value = 2044412.918677504
$ws.Cells.Item(9, 19).Value = 0.3413567538304502
$ws.Cells.Item(9, 20).Value = 0.3413567538304502

$ws.Cells.Item(10, 9).Value = 0.8830494168872806
$ws.Cells.Item(10, 10).Value = 0.8830494168872804
$ws.Cells.Item(10, 13).Value = 41.09915599999999
$ws.Cells.Item(10, 14).Value = 123.297468
$ws.Cells.Item(10, 15).Value = 0.09521587377309249
$ws.Cells.Item(10, 16).Value = 0.09521587377309249
$ws.Cells.Item(10, 17).Value = 55951.53073196575
$ws.Cells.Item(10, 18).Value = 503563.7765876917
$ws.Cells.Item(10, 19).Value = 0.08408032181374223
$ws.Cells.Item(10, 20).Value = 0.08408032181374221

$ws.Cells.Item(11, 9).Value = 0.8830494168872806
$ws.Cells.Item(11, 10).Value = 0.8830494168872804
$ws.Cells.Item(11, 13).Value = 112.0244103333333
$ws.Cells.Item(11, 14).Value = 336.073231
$ws.Cells.Item(11, 15).Value = 0.2595309284162377
$ws.Cells.Item(11, 16).Value = 0.2595309284162377
$ws.Cells.Item(11, 17).Value = 152507.6874448673
$ws.Cells.Item(11, 18).Value = 1372569.187003806
$ws.Cells.Item(11, 19).Value = 0.2291786350021733
$ws.Cells.Item(11, 20).Value = 0.2291786350021732

$ws.Cells.Item(12, 7).Value = 44.831112
$ws.Cells.Item(12, 8).Value = 134.493336
$ws.Cells.Item(12, 9).Value = 0.02907940059566787
$ws.Cells.Item(12, 10).Value = 0.02907940059566786
$ws.Cells.Item(12, 13).Value = 19.21315233333334
$ws.Cells.Item(12, 14).Value = 57.63945700000001
$ws.Cells.Item(12, 15).Value = 0.04451179209991234
$ws.Cells.Item(12, 16).Value = 0.04451179209991233
$ws.Cells.Item(12, 17).Value = 861.3469841287281
$ws.Cells.Item(12, 18).Value = 7752.122857158553
$ws.Cells.Item(12, 19).Value = 0.001294376233704435
$ws.Cells.Item(12, 20).Value = 0.001294376233704435

$ws.Cells.Item(13, 7).Value = 44.831112
$ws.Cells.Item(13, 8).Value = 134.493336
$ws.Cells.Item(13, 9).Value = 0.02907940059566787
$ws.Cells.Item(13, 10).Value = 0.02907940059566786
$ws.Cells.Item(13, 15).Value = 0.2141755495962477
$ws.Cells.Item(13, 16).Value = 0.2141755495962477
$ws.Cells.Item(13, 17).Value = 4144.507668995968
$ws.Cells.Item(13, 18).Value = 37300.56902096371
$ws.Cells.Item(13, 19).Value = 0.00622809660450662
$ws.Cells.Item(13, 20).Value = 0.006228096604506618

$ws.Cells.Item(14, 7).Value = 44.831112
$ws.Cells.Item(14, 8).Value = 134.493336
$ws.Cells.Item(14, 9).Value = 0.02907940059566787
$ws.Cells.Item(14, 10).Value = 0.02907940059566786
$ws.Cells.Item(14, 13).Value = 166.8580016666666
$ws.Cells.Item(14, 14).Value = 500.5740049999999
$ws.Cells.Item(14, 15).Value = 0.3865658561145097
$ws.Cells.Item(14, 16).Value = 0.3865658561145097
$ws.Cells.Item(14, 17).Value = 7480.429760814518
$ws.Cells.Item(14, 18).Value = 67323.86784733068
$ws.Cells.Item(14, 19).Value = 0.01124110338656113
$ws.Cells.Item(14, 20).Value = 0.01124110338656113

$ws.Cells.Item(15, 7).Value = 44.831112
$ws.Cells.Item(15, 8).Value = 134.493336
$ws.Cells.Item(15, 9).Value = 0.02907940059566787
$ws.Cells.Item(15, 10).Value = 0.02907940059566786
$ws.Cells.Item(15, 13).Value = 41.09915599999999
$ws.Cells.Item(15, 14).Value = 123.297468
$ws.Cells.Item(15, 15).Value = 0.09521587377309249
$ws.Cells.Item(15, 16).Value = 0.09521587377309249
$ws.Cells.Item(15, 17).Value = 1842.520865741472
$ws.Cells.Item(15, 18).Value = 16582.68779167324
$ws.Cells.Item(15, 19).Value = 0.002768820536514303
$ws.Cells.Item(15, 20).Value = 0.002768820536514302

$ws.Cells.Item(16, 7).Value = 44.831112
$ws.Cells.Item(16, 8).Value = 134.493336
$ws.Cells.Item(16, 9).Value = 0.02907940059566787
$ws.Cells.Item(16, 10).Value = 0.02907940059566786
$ws.Cells.Item(16, 13).Value = 112.0244103333333
$ws.Cells.Item(16, 14).Value = 336.073231
$ws.Cells.Item(16, 15).Value = 0.2595309284162377
$ws.Cells.Item(16, 16).Value = 0.2595309284162377
$ws.Cells.Item(16, 17).Value = 5022.178886387624
$ws.Cells.Item(16, 18).Value = 45199.60997748862
$ws.Cells.Item(16, 19).Value = 0.007547003834381378
$ws.Cells.Item(16, 20).Value = 0.007547003834381377

$ws.Cells.Item(17, 7).Value = 52.83062100000001
$ws.Cells.Item(17, 8).Value = 158.491863
$ws.Cells.Item(17, 9).Value = 0.0342682285413064
$ws.Cells.Item(17, 10).Value = 0.03426822854130639
$ws.Cells.Item(17, 13).Value = 19.21315233333334
$ws.Cells.Item(17, 14).Value = 57.63945700000001
$ws.Cells.Item(17, 15).Value = 0.04451179209991234
$ws.Cells.Item(17, 16).Value = 0.04451179209991233
$ws.Cells.Item(17, 17).Value = 1015.042769137599
$ws.Cells.Item(17, 18).Value = 9135.384922238394
$ws.Cells.Item(17, 19).Value = 0.001525340264462913
$ws.Cells.Item(17, 20).Value = 0.001525340264462912

$ws.Cells.Item(18, 7).Value = 52.83062100000001
$ws.Cells.Item(18, 8).Value = 158.491863
$ws.Cells.Item(18, 9).Value = 0.0342682285413064
$ws.Cells.Item(18, 10).Value = 0.03426822854130639
$ws.Cells.Item(18, 15).Value = 0.2141755495962477
$ws.Cells.Item(18, 16).Value = 0.2141755495962477
$ws.Cells.Item(18, 17).Value = 4884.039322788145
$ws.Cells.Item(18, 18).Value = 43956.3539050933
$ws.Cells.Item(18, 19).Value = 0.007339416681524122
$ws.Cells.Item(18, 20).Value = 0.007339416681524119

$ws.Cells.Item(19, 7).Value = 52.83062100000001
$ws.Cells.Item(19, 8).Value = 158.491863
$ws.Cells.Item(19, 9).Value = 0.0342682285413064
$ws.Cells.Item(19, 10).Value = 0.03426822854130639
$ws.Cells.Item(19, 13).Value = 166.8580016666666
$ws.Cells.Item(19, 14).Value = 500.5740049999999
$ws.Cells.Item(19, 15).Value = 0.3865658561145097
$ws.Cells.Item(19, 16).Value = 0.3865658561145097
$ws.Cells.Item(19, 17).Value = 8815.211846869035
$ws.Cells.Item(19, 18).Value = 79336.90662182131
$ws.Cells.Item(19, 19).Value = 0.01324692710359779
$ws.Cells.Item(19, 20).Value = 0.01324692710359778

$ws.Cells.Item(20, 7).Value = 52.83062100000001
$ws.Cells.Item(20, 8).Value = 158.491863
$ws.Cells.Item(20, 9).Value = 0.0342682285413064
$ws.Cells.Item(20, 10).Value = 0.03426822854130639
$ws.Cells.Item(20, 13).Value = 41.09915599999999
$ws.Cells.Item(20, 14).Value = 123.297468
$ws.Cells.Item(20, 15).Value = 0.09521587377309249
$ws.Cells.Item(20, 16).Value = 0.09521587377309249
$ws.Cells.Item(20, 17).Value = 2171.293934055876
$ws.Cells.Item(20, 18).Value = 19541.64540650288
$ws.Cells.Item(20, 19).Value = 0.003262879323216515
$ws.Cells.Item(20, 20).Value = 0.003262879323216515

$ws.Cells.Item(21, 7).Value = 52.83062100000001
$ws.Cells.Item(21, 8).Value = 158.491863
$ws.Cells.Item(21, 9).Value = 0.0342682285413064
$ws.Cells.Item(21, 10).Value = 0.03426822854130639
$ws.Cells.Item(21, 13).Value = 112.0244103333333
$ws.Cells.Item(21, 14).Value = 336.073231
$ws.Cells.Item(21, 15).Value = 0.2595309284162377
$ws.Cells.Item(21, 16).Value = 0.2595309284162377
$ws.Cells.Item(21, 17).Value = 5918.319165068818
$ws.Cells.Item(21, 18).Value = 53264.87248561937
$ws.Cells.Item(21, 19).Value = 0.008893665168505066
$ws.Cells.Item(21, 20).Value = 0.008893665168505063

$ws.Cells.Item(22, 7).Value = 16.16161433333333
$ws.Cells.Item(22, 8).Value = 48.484843
$ws.Cells.Item(22, 9).Value = 0.01048312291409786
$ws.Cells.Item(22, 10).Value = 0.01048312291409786
$ws.Cells.Item(22, 13).Value = 19.21315233333334
$ws.Cells.Item(22, 14).Value = 57.63945700000001
$ws.Cells.Item(22, 15).Value = 0.04451179209991234
$ws.Cells.Item(22, 16).Value = 0.04451179209991233
$ws.Cells.Item(22, 17).Value = 310.5155581389168
$ws.Cells.Item(22, 18).Value = 2794.640023250251
$ws.Cells.Item(22, 19).Value = 0.0004666225877101513
$ws.Cells.Item(22, 20).Value = 0.0004666225877101511

$ws.Cells.Item(23, 7).Value = 16.16161433333333
$ws.Cells.Item(23, 8).Value = 48.484843
$ws.Cells.Item(23, 9).Value = 0.01048312291409786
$ws.Cells.Item(23, 10).Value = 0.01048312291409786
$ws.Cells.Item(23, 15).Value = 0.2141755495962477
$ws.Cells.Item(23, 16).Value = 0.2141755495962477
$ws.Cells.Item(23, 17).Value = 1494.094872057939
$ws.Cells.Item(23, 18).Value = 13446.85384852145
$ws.Cells.Item(23, 19).Value = 0.002245228611611928
$ws.Cells.Item(23, 20).Value = 0.002245228611611927

$ws.Cells.Item(24, 7).Value = 16.16161433333333
$ws.Cells.Item(24, 8).Value = 48.484843
$ws.Cells.Item(24, 9).Value = 0.01048312291409786
$ws.Cells.Item(24, 10).Value = 0.01048312291409786
$ws.Cells.Item(24, 13).Value = 166.8580016666666
$ws.Cells.Item(24, 14).Value = 500.5740049999999
$ws.Cells.Item(24, 15).Value = 0.3865658561145097
$ws.Cells.Item(24, 16).Value = 0.3865658561145097
$ws.Cells.Item(24, 17).Value = 2696.694671367356
$ws.Cells.Item(24, 18).Value = 24270.25204230621
$ws.Cells.Item(24, 19).Value = 0.004052417384041875
$ws.Cells.Item(24, 20).Value = 0.004052417384041874

$ws.Cells.Item(25, 7).Value = 16.16161433333333
$ws.Cells.Item(25, 8).Value = 48.484843
$ws.Cells.Item(25, 9).Value = 0.01048312291409786
$ws.Cells.Item(25, 10).Value = 0.01048312291409786
$ws.Cells.Item(25, 13).Value = 41.09915599999999
$ws.Cells.Item(25, 14).Value = 123.297468
$ws.Cells.Item(25, 15).Value = 0.09521587377309249
$ws.Cells.Item(25, 16).Value = 0.09521587377309249
$ws.Cells.Item(25, 17).Value = 664.2287086975025
$ws.Cells.Item(25, 18).Value = 5978.058378277523
$ws.Cells.Item(25, 19).Value = 0.0009981597081365557
$ws.Cells.Item(25, 20).Value = 0.0009981597081365555

$ws.Cells.Item(26, 7).Value = 16.16161433333333
$ws.Cells.Item(26, 8).Value = 48.484843
$ws.Cells.Item(26, 9).Value = 0.01048312291409786
$ws.Cells.Item(26, 10).Value = 0.01048312291409786
$ws.Cells.Item(26, 13).Value = 112.0244103333333
$ws.Cells.Item(26, 14).Value = 336.073231
$ws.Cells.Item(26, 15).Value = 0.2595309284162377
$ws.Cells.Item(26, 16).Value = 0.2595309284162377
$ws.Cells.Item(26, 17).Value = 1810.495315726415
$ws.Cells.Item(26, 18).Value = 16294.45784153773
$ws.Cells.Item(26, 19).Value = 0.002720694622597354
$ws.Cells.Item(26, 20).Value = 0.002720694622597354
